$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-12 Saturday" "2024-10-13 Sunday"

Replace-Text "64×47=" "23×12="
Replace-Text "71×84=" "89×94="
Replace-Text "83×58=" "48×23="
Replace-Text "65×13=" "46×96="
Replace-Text "94×18=" "82×89="

Replace-Text "29×16=" "67×17="
Replace-Text "86×24=" "98×27="
Replace-Text "95×78=" "66×33="
Replace-Text "18×79=" "21×30="
Replace-Text "22×80=" "22×27="

Replace-Text "69×23=" "80×81="
Replace-Text "57×34=" "46×69="
Replace-Text "55×59=" "45×84="
Replace-Text "38×67=" "91×28="
Replace-Text "45×83=" "36×11="

Replace-Text "51×21=" "96×60="
Replace-Text "14×81=" "37×78="
Replace-Text "87×77=" "32×74="
Replace-Text "93×76=" "24×67="
Replace-Text "89×25=" "14×37="

Replace-Text "47×74=" "98×82="
Replace-Text "96×95=" "94×34="
Replace-Text "40×82=" "60×81="
Replace-Text "72×91=" "92×54="
Replace-Text "59×93=" "36×18="
